$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append new row 51 with the second test mail ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(51, 1).Value = "Wil je 100 stuks M5-bouten bestellen?"
$logs.Cells.Item(51, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(51, 3).Value = "Testmail #2: Wil je 100 stuks M5-bouten bestellen?"
$logs.Cells.Item(51, 4).Value = "Inkoop / Bestellingen"
$logs.Cells.Item(51, 5).Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Cells.Item(51, 6).Value = "2025-08-03 23:28:53"
$logs.Cells.Item(51, 7).Value = "Ja"
$logs.Cells.Item(51, 8).Value = "Ja"
$logs.Cells.Item(51, 9).Value = "Nee"
$logs.Cells.Item(51, 10).Value = "Nee"

# --- Extend the conditional formatting ranges on Logs from row 50 to row 51 ---
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$col`2:$col`50")
    $newRange = $logs.Range("$col`2:$col`51")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- "Dashboard" sheet: bump the "Inkoop / Bestellingen" count from 9 to 10 ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(4, 2).Value = 10
